# Auto-generated Excel COM-interop script
# Applies numeric cell updates (and a few cell clears) across 8 worksheets
# to match the target diff for Balmung_Profits market data.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

$updates = @(
    @{ Cell = "H9"; Value = 5558848.5 }
    @{ Cell = "I9"; Value = 6410819.5 }
    @{ Cell = "K9"; Value = 6410819.5 }
    @{ Cell = "M9"; Value = -6410650.5 }
    @{ Cell = "H135"; Value = 1172.1666 }
    @{ Cell = "I135"; Value = 1172.1666 }
    @{ Cell = "K135"; Value = 10549.4994 }
    @{ Cell = "M135"; Value = -8014.499400000001 }
    @{ Cell = "H137"; Value = 1322139.8 }
    @{ Cell = "I137"; Value = 6243.9116 }
    @{ Cell = "J137"; Value = 2387388.8 }
    @{ Cell = "K137"; Value = 18731.7348 }
    @{ Cell = "L137"; Value = 7162166.399999999 }
    @{ Cell = "M137"; Value = -16181.7348 }
    @{ Cell = "N137"; Value = -7167266.399999999 }
    @{ Cell = "H138"; Value = 20131.572 }
    @{ Cell = "J138"; Value = 4584.2 }
    @{ Cell = "L138"; Value = 13752.6 }
    @{ Cell = "N138"; Value = -24032.6 }
)
foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

$updates = @(
    @{ Cell = "H32"; Value = 2222.4482 }
    @{ Cell = "I32"; Value = 1171.8784 }
    @{ Cell = "K32"; Value = 1171.8784 }
    @{ Cell = "M32"; Value = -884.8784000000001 }
    @{ Cell = "H37"; Value = 23999.8 }
    @{ Cell = "I37"; Value = 14999.5 }
    @{ Cell = "J37"; Value = 30000 }
    @{ Cell = "K37"; Value = 14999.5 }
    @{ Cell = "L37"; Value = 30000 }
    @{ Cell = "M37"; Value = -14726.5 }
    @{ Cell = "N37"; Value = -30546 }
    @{ Cell = "H55"; Value = 5000 }
    @{ Cell = "I55"; Value = 5000 }
    @{ Cell = "K55"; Value = 5000 }
    @{ Cell = "M55"; Value = -4685 }
    @{ Cell = "H61"; Value = 1429678.1 }
    @{ Cell = "I61"; Value = 38031.867 }
    @{ Cell = "J61"; Value = 4908794 }
    @{ Cell = "K61"; Value = 38031.867 }
    @{ Cell = "L61"; Value = 4908794 }
    @{ Cell = "M61"; Value = -37819.867 }
    @{ Cell = "N61"; Value = -4909218 }
    @{ Cell = "H74"; Value = 657696.3 }
    @{ Cell = "I74"; Value = 3751.9285 }
    @{ Cell = "J74"; Value = 3709436.8 }
    @{ Cell = "K74"; Value = 3751.9285 }
    @{ Cell = "L74"; Value = 3709436.8 }
    @{ Cell = "M74"; Value = -2877.9285 }
    @{ Cell = "N74"; Value = -3711184.8 }
    @{ Cell = "H77"; Value = 657696.3 }
    @{ Cell = "I77"; Value = 3751.9285 }
    @{ Cell = "J77"; Value = 3709436.8 }
    @{ Cell = "K77"; Value = 18759.6425 }
    @{ Cell = "L77"; Value = 18547184 }
    @{ Cell = "M77"; Value = -14391.6425 }
    @{ Cell = "N77"; Value = -18555920 }
    @{ Cell = "H80"; Value = 0 }
    @{ Cell = "J80"; Value = 0 }
    @{ Cell = "L80"; Value = 0 }
    @{ Cell = "H83"; Value = 0 }
    @{ Cell = "J83"; Value = 0 }
    @{ Cell = "L83"; Value = 0 }
    @{ Cell = "H88"; Value = 2847 }
    @{ Cell = "I88"; Value = 0 }
    @{ Cell = "J88"; Value = 2847 }
    @{ Cell = "K88"; Value = 0 }
    @{ Cell = "L88"; Value = 2847 }
    @{ Cell = "N88"; Value = -3659 }
    @{ Cell = "H91"; Value = 2847 }
    @{ Cell = "I91"; Value = 0 }
    @{ Cell = "J91"; Value = 2847 }
    @{ Cell = "K91"; Value = 0 }
    @{ Cell = "L91"; Value = 2847 }
    @{ Cell = "N91"; Value = -5655 }
    @{ Cell = "H105"; Value = 97250 }
    @{ Cell = "I105"; Value = 97250 }
    @{ Cell = "K105"; Value = 97250 }
    @{ Cell = "M105"; Value = -93756 }
    @{ Cell = "H136"; Value = 1429678.1 }
    @{ Cell = "I136"; Value = 38031.867 }
    @{ Cell = "J136"; Value = 4908794 }
    @{ Cell = "K136"; Value = 114095.601 }
    @{ Cell = "L136"; Value = 14726382 }
    @{ Cell = "M136"; Value = -111545.601 }
    @{ Cell = "N136"; Value = -14731482 }
)
foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

$clears = @("N80", "N83", "M88", "M91")
foreach ($r in $clears) {
    $ws.Range($r).ClearContents()
}

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

$updates = @(
    @{ Cell = "H86"; Value = 4808.826 }
    @{ Cell = "I86"; Value = 3153.2354 }
    @{ Cell = "J86"; Value = 9499.666999999999 }
    @{ Cell = "K86"; Value = 3153.2354 }
    @{ Cell = "L86"; Value = 9499.666999999999 }
    @{ Cell = "M86"; Value = -2030.2354 }
    @{ Cell = "N86"; Value = -11745.667 }
    @{ Cell = "H89"; Value = 4808.826 }
    @{ Cell = "I89"; Value = 3153.2354 }
    @{ Cell = "J89"; Value = 9499.666999999999 }
    @{ Cell = "K89"; Value = 15766.177 }
    @{ Cell = "L89"; Value = 47498.335 }
    @{ Cell = "M89"; Value = -10150.177 }
    @{ Cell = "N89"; Value = -58730.335 }
    @{ Cell = "H107"; Value = 20343.53 }
    @{ Cell = "I107"; Value = 22842.615 }
    @{ Cell = "K107"; Value = 22842.615 }
    @{ Cell = "M107"; Value = -20922.615 }
    @{ Cell = "H109"; Value = 174777 }
    @{ Cell = "J109"; Value = 174777 }
    @{ Cell = "L109"; Value = 174777 }
    @{ Cell = "N109"; Value = -177551 }
)
foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

$updates = @(
    @{ Cell = "H22"; Value = 611.56525 }
    @{ Cell = "I22"; Value = 626.9545000000001 }
    @{ Cell = "K22"; Value = 626.9545000000001 }
    @{ Cell = "M22"; Value = -276.9545000000001 }
    @{ Cell = "H31"; Value = 4573.3896 }
    @{ Cell = "I31"; Value = 2774.9 }
    @{ Cell = "J31"; Value = 5052.987 }
    @{ Cell = "K31"; Value = 2774.9 }
    @{ Cell = "L31"; Value = 5052.987 }
    @{ Cell = "M31"; Value = -2479.9 }
    @{ Cell = "N31"; Value = -5642.987 }
    @{ Cell = "H34"; Value = 4573.3896 }
    @{ Cell = "I34"; Value = 2774.9 }
    @{ Cell = "J34"; Value = 5052.987 }
    @{ Cell = "K34"; Value = 2774.9 }
    @{ Cell = "L34"; Value = 5052.987 }
    @{ Cell = "M34"; Value = -2572.9 }
    @{ Cell = "N34"; Value = -5456.987 }
    @{ Cell = "H107"; Value = 1744.6666 }
    @{ Cell = "J107"; Value = 2199 }
    @{ Cell = "L107"; Value = 2199 }
    @{ Cell = "N107"; Value = -6039 }
    @{ Cell = "H108"; Value = 0 }
    @{ Cell = "I108"; Value = 0 }
    @{ Cell = "K108"; Value = 0 }
    @{ Cell = "H109"; Value = 49999.75 }
    @{ Cell = "J109"; Value = 49999.75 }
    @{ Cell = "L109"; Value = 49999.75 }
    @{ Cell = "N109"; Value = -52079.75 }
    @{ Cell = "H132"; Value = 2763.64 }
    @{ Cell = "I132"; Value = 2927.4443 }
    @{ Cell = "J132"; Value = 2342.4285 }
    @{ Cell = "K132"; Value = 8782.332900000001 }
    @{ Cell = "L132"; Value = 7027.2855 }
    @{ Cell = "M132"; Value = -6252.332900000001 }
    @{ Cell = "N132"; Value = -12087.2855 }
    @{ Cell = "H138"; Value = 93615.17999999999 }
    @{ Cell = "J138"; Value = 93615.17999999999 }
    @{ Cell = "L138"; Value = 93615.17999999999 }
    @{ Cell = "N138"; Value = -103895.18 }
)
foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

$clears = @("M108")
foreach ($r in $clears) {
    $ws.Range($r).ClearContents()
}

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

$updates = @(
    @{ Cell = "H32"; Value = 14829.167 }
    @{ Cell = "I32"; Value = 18982 }
    @{ Cell = "J32"; Value = 13998.6 }
    @{ Cell = "K32"; Value = 56946 }
    @{ Cell = "L32"; Value = 41995.8 }
    @{ Cell = "M32"; Value = -56663 }
    @{ Cell = "N32"; Value = -42561.8 }
    @{ Cell = "H46"; Value = 11111833 }
    @{ Cell = "I46"; Value = 100000000 }
    @{ Cell = "J46"; Value = 812.375 }
    @{ Cell = "K46"; Value = 300000000 }
    @{ Cell = "L46"; Value = 2437.125 }
    @{ Cell = "M46"; Value = -299999909 }
    @{ Cell = "N46"; Value = -2619.125 }
    @{ Cell = "H119"; Value = 11041.917 }
    @{ Cell = "I119"; Value = 4564.75 }
    @{ Cell = "K119"; Value = 13694.25 }
    @{ Cell = "M119"; Value = -8856.25 }
    @{ Cell = "H121"; Value = 71441700 }
    @{ Cell = "I121"; Value = 250001520 }
    @{ Cell = "J121"; Value = 17766.2 }
    @{ Cell = "K121"; Value = 750004560 }
    @{ Cell = "L121"; Value = 53298.60000000001 }
    @{ Cell = "M121"; Value = -750003250 }
    @{ Cell = "N121"; Value = -55918.60000000001 }
    @{ Cell = "H122"; Value = 10775755 }
    @{ Cell = "I122"; Value = 12122049 }
    @{ Cell = "J122"; Value = 5405 }
    @{ Cell = "K122"; Value = 109098441 }
    @{ Cell = "L122"; Value = 48645 }
    @{ Cell = "M122"; Value = -109095991 }
    @{ Cell = "N122"; Value = -53545 }
    @{ Cell = "H123"; Value = 8166.6665 }
    @{ Cell = "I123"; Value = 4500 }
    @{ Cell = "J123"; Value = 10000 }
    @{ Cell = "K123"; Value = 13500 }
    @{ Cell = "L123"; Value = 30000 }
    @{ Cell = "M123"; Value = -11050 }
    @{ Cell = "N123"; Value = -34900 }
)
foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

$updates = @(
    @{ Cell = "H105"; Value = 55113 }
    @{ Cell = "J105"; Value = 55113 }
    @{ Cell = "L105"; Value = 55113 }
    @{ Cell = "N105"; Value = -62101 }
    @{ Cell = "H107"; Value = 168564.5 }
    @{ Cell = "I107"; Value = 250346.75 }
    @{ Cell = "K107"; Value = 250346.75 }
    @{ Cell = "M107"; Value = -248426.75 }
    @{ Cell = "H108"; Value = 158888.5 }
    @{ Cell = "J108"; Value = 158888.5 }
    @{ Cell = "L108"; Value = 158888.5 }
    @{ Cell = "N108"; Value = -166568.5 }
)
foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

$updates = @(
    @{ Cell = "H46"; Value = 6538.15 }
    @{ Cell = "I46"; Value = 27449.5 }
    @{ Cell = "J46"; Value = 1310.3125 }
    @{ Cell = "K46"; Value = 27449.5 }
    @{ Cell = "L46"; Value = 1310.3125 }
    @{ Cell = "M46"; Value = -27261.5 }
    @{ Cell = "N46"; Value = -1686.3125 }
    @{ Cell = "H55"; Value = 704.80646 }
    @{ Cell = "I55"; Value = 515.73334 }
    @{ Cell = "J55"; Value = 882.0625 }
    @{ Cell = "K55"; Value = 515.73334 }
    @{ Cell = "L55"; Value = 882.0625 }
    @{ Cell = "M55"; Value = -342.73334 }
    @{ Cell = "N55"; Value = -1228.0625 }
    @{ Cell = "H106"; Value = 16621.8 }
    @{ Cell = "J106"; Value = 16621.8 }
    @{ Cell = "L106"; Value = 16621.8 }
    @{ Cell = "N106"; Value = -19145.8 }
    @{ Cell = "H109"; Value = 87083.25 }
    @{ Cell = "J109"; Value = 87083.25 }
    @{ Cell = "L109"; Value = 87083.25 }
    @{ Cell = "N109"; Value = -89857.25 }
    @{ Cell = "H122"; Value = 3229.389 }
    @{ Cell = "I122"; Value = 2843.76 }
    @{ Cell = "J122"; Value = 4105.8184 }
    @{ Cell = "K122"; Value = 8531.280000000001 }
    @{ Cell = "L122"; Value = 12317.4552 }
    @{ Cell = "M122"; Value = -6081.280000000001 }
    @{ Cell = "N122"; Value = -17217.4552 }
    @{ Cell = "H132"; Value = 4327.2334 }
    @{ Cell = "I132"; Value = 3820.92 }
    @{ Cell = "J132"; Value = 6858.8 }
    @{ Cell = "K132"; Value = 11462.76 }
    @{ Cell = "L132"; Value = 20576.4 }
    @{ Cell = "M132"; Value = -8932.76 }
    @{ Cell = "N132"; Value = -25636.4 }
)
foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

$updates = @(
    @{ Cell = "H105"; Value = 44999.668 }
    @{ Cell = "I105"; Value = 55000 }
    @{ Cell = "J105"; Value = 39999.5 }
    @{ Cell = "K105"; Value = 55000 }
    @{ Cell = "L105"; Value = 39999.5 }
    @{ Cell = "M105"; Value = -51506 }
    @{ Cell = "N105"; Value = -46987.5 }
    @{ Cell = "H126"; Value = 2026.7368 }
    @{ Cell = "I126"; Value = 2083.7778 }
    @{ Cell = "J126"; Value = 1000 }
    @{ Cell = "K126"; Value = 6251.3334 }
    @{ Cell = "L126"; Value = 3000 }
    @{ Cell = "M126"; Value = -3781.3334 }
    @{ Cell = "N126"; Value = -7940 }
    @{ Cell = "H132"; Value = 2842.75 }
    @{ Cell = "I132"; Value = 2417.5 }
    @{ Cell = "J132"; Value = 4969 }
    @{ Cell = "K132"; Value = 7252.5 }
    @{ Cell = "L132"; Value = 14907 }
    @{ Cell = "M132"; Value = -4722.5 }
    @{ Cell = "N132"; Value = -19967 }
)
foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
